$wb = $excel.ActiveWorkbook

# Rename Sheet2 -> LoginData
$wsLogin = $wb.Worksheets.Item("Sheet2")
$wsLogin.Name = "LoginData"

# Delete the two leading blank rows (1:2) on LoginData so data shifts from rows 3-10 to rows 1-8
$wsLogin.Rows("1:2").Delete()
